# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new survey week) above the existing
# row 15, pushing the former rows 15-38 down to 17-40, then fill in the
# new rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 15 (existing rows 15-38 shift down to 17-40)
$ws.Rows("15:16").Insert()

# New row 15
$ws.Cells(15, 1).Value = 10
$ws.Cells(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells(15, 3).Value = "La Araucanía"
$ws.Cells(15, 4).Value = 45260
$ws.Cells(15, 5).Value = 9
$ws.Cells(15, 6).Value = "Fruta"
$ws.Cells(15, 7).Value = 100104
$ws.Cells(15, 8).Value = "Frutos de pepita"
$ws.Cells(15, 9).Value = 100104004
$ws.Cells(15, 10).Value = "Níspero"
$ws.Cells(15, 11).Value = "Californiana(o)"
$ws.Cells(15, 12).Value = "Primera"
$ws.Cells(15, 13).Value = 180
$ws.Cells(15, 14).Value = 28000
$ws.Cells(15, 15).Value = 28000
$ws.Cells(15, 16).Value = 28000
$ws.Cells(15, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells(15, 18).Value = "Provincia de Quillota"
$ws.Cells(15, 19).Value = 2800
$ws.Cells(15, 20).Value = 10

# New row 16
$ws.Cells(16, 1).Value = 10
$ws.Cells(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells(16, 3).Value = "La Araucanía"
$ws.Cells(16, 4).Value = 45260
$ws.Cells(16, 5).Value = 9
$ws.Cells(16, 6).Value = "Fruta"
$ws.Cells(16, 7).Value = 100104
$ws.Cells(16, 8).Value = "Frutos de pepita"
$ws.Cells(16, 9).Value = 100104004
$ws.Cells(16, 10).Value = "Níspero"
$ws.Cells(16, 11).Value = "Californiana(o)"
$ws.Cells(16, 12).Value = "Primera"
$ws.Cells(16, 13).Value = 400
$ws.Cells(16, 14).Value = 20000
$ws.Cells(16, 15).Value = 20000
$ws.Cells(16, 16).Value = 20000
$ws.Cells(16, 17).Value = "`$/bandeja 5 kilos"
$ws.Cells(16, 18).Value = "Provincia de Quillota"
$ws.Cells(16, 19).Value = 4000
$ws.Cells(16, 20).Value = 5
